$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 42 (new betting entry: WWE entertainment) ---
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 33
$ws.Range("C42").Value = 45318
$ws.Range("D42").Value = 1.08
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 778
$ws.Range("G42").Formula = "=F42*E42*D42"
$ws.Range("H42").Formula = "=I41"
$ws.Range("I42").Formula = "=H42+G42-F42"
$ws.Range("J42").Formula = "=I42-H42"
$ws.Range("K42").Formula = '=I42/$H$2-1'
$ws.Range("L42").Value = "NA"
$ws.Range("M42").Value = "ENTERTAINMENT"
$ws.Range("N42").Value = "WWE"

# --- Row 43 (new betting entry: WWE entertainment) ---
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 34
$ws.Range("C43").Value = 45318
$ws.Range("D43").Value = 1.1
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = 6000
$ws.Range("G43").Formula = "=F43*E43*D43"
$ws.Range("H43").Formula = "=I42"
$ws.Range("I43").Formula = "=H43+G43-F43"
$ws.Range("J43").Formula = "=I43-H43"
$ws.Range("K43").Formula = '=I43/$H$2-1'
$ws.Range("L43").Value = "NA"
$ws.Range("M43").Value = "ENTERTAINMENT"
$ws.Range("N43").Value = "WWE"

# --- Clone number formats from the row above (reuses existing style records) ---
$ws.Range("C41:D41").Copy()
$ws.Range("C42:D43").PasteSpecial(-4122)
$ws.Range("G41").Copy()
$ws.Range("G42:G43").PasteSpecial(-4122)
$ws.Range("H41:K41").Copy()
$ws.Range("H42:K43").PasteSpecial(-4122)

# --- Column M widens to fit "ENTERTAINMENT" ---
$ws.Columns.Item(13).ColumnWidth = 14.3

# --- View / selection state to match the saved workbook ---
$excel.ActiveWindow.ScrollRow = 36
$ws.Range("F44").Select()
